$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Commit-ID (column F) back-fill for existing rows ---
$ws.Range("F17").Value = "648f0de37aa2386fc4ae2455ef79559372560b9b"
$ws.Range("F30").Value = "1f664c10de6556ca5e4d14a3e5c612d552ac1ba2"
$ws.Range("F43").Value = "be538551470e444abf2b2d3ed25d7d6d762b2ee1"
$ws.Range("F53").Value = "b309fe8dd426613962ad859ff269a4a711c37c04"

# --- New rows 54-57 : work diary entries for 16.05.2022 (serial 44697) ---
# Carry the date-number format down from the last existing row, then fill values.
$ws.Range("A53").Copy()
$ws.Range("A54:A57").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newDate = 44697
$ws.Range("A54").Value = $newDate
$ws.Range("A55").Value = $newDate
$ws.Range("A56").Value = $newDate
$ws.Range("A57").Value = $newDate

# Row 54
$ws.Range("B54").Value = "Analyse"
$ws.Range("C54").Value = 1
$ws.Range("D54").Value = "Documentation du projet"

# Row 55
$ws.Range("B55").Value = "Réalisation"
$ws.Range("C55").Value = 0.5
$ws.Range("D55").Value = "ajout de tâches dans le scrum, définitions de tests"

# Row 56
$ws.Range("B56").Value = "Conception"
$ws.Range("C56").Value = 2.25
$ws.Range("D56").Value = "Création des maquettes pour ce sprint"
$ws.Range("F56").Value = "65f4ed811c4dc8076e48ec6f7af41212b192fd81"

# Row 57
$ws.Range("B57").Value = "analyse"
$ws.Range("C57").Value = 1.5
$ws.Range("D57").Value = "documentation"

# --- Grow the table / autofilter to cover the new rows ---
$table = $ws.ListObjects.Item("Tableau1")
$table.Resize($ws.Range("A1:F57"))

# --- Update selection to match the edited area ---
$ws.Range("E57").Select()
